$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 42.714287
$ws.Range("I8").Value = 30.5
$ws.Range("J8").Value = 59
$ws.Range("K8").Value = 91.5
$ws.Range("L8").Value = 177
$ws.Range("M8").Value = 47.5
$ws.Range("N8").Value = -455
$ws.Range("H40").Value = 62502436
$ws.Range("J40").Value = 62502436
$ws.Range("L40").Value = 62502436
$ws.Range("N40").Value = -62502786
$ws.Range("H43").Value = 4666
$ws.Range("I43").Value = 4499
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 4499
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -4430
$ws.Range("N43").Value = -5138
$ws.Range("H55").Value = 227.1875
$ws.Range("I55").Value = 174.33333
$ws.Range("J55").Value = 258.9
$ws.Range("K55").Value = 174.33333
$ws.Range("L55").Value = 258.9
$ws.Range("M55").Value = 39.66667000000001
$ws.Range("N55").Value = -686.9
$ws.Range("H86").Value = 1957.7858
$ws.Range("I86").Value = 1375.125
$ws.Range("J86").Value = 2734.6667
$ws.Range("K86").Value = 1375.125
$ws.Range("L86").Value = 2734.6667
$ws.Range("M86").Value = -252.125
$ws.Range("N86").Value = -4980.6667
$ws.Range("H89").Value = 1957.7858
$ws.Range("I89").Value = 1375.125
$ws.Range("J89").Value = 2734.6667
$ws.Range("K89").Value = 6875.625
$ws.Range("L89").Value = 13673.3335
$ws.Range("M89").Value = -1259.625
$ws.Range("N89").Value = -24905.3335
$ws.Range("H116").Value = 18074.25
$ws.Range("I116").Value = 14099
$ws.Range("J116").Value = 30000
$ws.Range("K116").Value = 14099
$ws.Range("L116").Value = 30000
$ws.Range("M116").Value = -10657
$ws.Range("N116").Value = -36884

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 693.06665
$ws.Range("I2").Value = 381.42856
$ws.Range("J2").Value = 965.75
$ws.Range("K2").Value = 381.42856
$ws.Range("L2").Value = 965.75
$ws.Range("M2").Value = -268.42856
$ws.Range("N2").Value = -1191.75
$ws.Range("H16").Value = 1351.5
$ws.Range("J16").Value = 1633.3334
$ws.Range("L16").Value = 1633.3334
$ws.Range("N16").Value = -2207.3334
$ws.Range("H32").Value = 5707.143
$ws.Range("I32").Value = 4270.328
$ws.Range("K32").Value = 4270.328
$ws.Range("M32").Value = -3983.328
$ws.Range("H46").Value = 35992.832
$ws.Range("J46").Value = 35992.832
$ws.Range("L46").Value = 35992.832
$ws.Range("N46").Value = -36630.832
$ws.Range("H61").Value = 7509290.5
$ws.Range("I61").Value = 8009793.5
$ws.Range("K61").Value = 8009793.5
$ws.Range("M61").Value = -8009581.5
$ws.Range("H74").Value = 2391.8333
$ws.Range("I74").Value = 1989.7778
$ws.Range("J74").Value = 3598
$ws.Range("K74").Value = 1989.7778
$ws.Range("L74").Value = 3598
$ws.Range("M74").Value = -1115.7778
$ws.Range("N74").Value = -5346
$ws.Range("H77").Value = 2391.8333
$ws.Range("I77").Value = 1989.7778
$ws.Range("J77").Value = 3598
$ws.Range("K77").Value = 9948.889000000001
$ws.Range("L77").Value = 17990
$ws.Range("M77").Value = -5580.889000000001
$ws.Range("N77").Value = -26726
$ws.Range("H88").Value = 3127.6365
$ws.Range("I88").Value = 2566.6667
$ws.Range("K88").Value = 2566.6667
$ws.Range("M88").Value = -2160.6667
$ws.Range("H91").Value = 3127.6365
$ws.Range("I91").Value = 2566.6667
$ws.Range("K91").Value = 2566.6667
$ws.Range("M91").Value = -1162.6667
$ws.Range("H110").Value = 5207.2383
$ws.Range("I110").Value = 5725.7144
$ws.Range("K110").Value = 5725.7144
$ws.Range("M110").Value = -3680.7144
$ws.Range("H116").Value = 693.06665
$ws.Range("I116").Value = 381.42856
$ws.Range("J116").Value = 965.75
$ws.Range("K116").Value = 381.42856
$ws.Range("L116").Value = 965.75
$ws.Range("M116").Value = 1912.57144
$ws.Range("N116").Value = -5553.75
$ws.Range("H136").Value = 7509290.5
$ws.Range("I136").Value = 8009793.5
$ws.Range("K136").Value = 24029380.5
$ws.Range("M136").Value = -24026830.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 693.06665
$ws.Range("I3").Value = 381.42856
$ws.Range("J3").Value = 965.75
$ws.Range("K3").Value = 381.42856
$ws.Range("L3").Value = 965.75
$ws.Range("M3").Value = -267.42856
$ws.Range("N3").Value = -1193.75
$ws.Range("H20").Value = 4945.1
$ws.Range("I20").Value = 4827.8887
$ws.Range("J20").Value = 6000
$ws.Range("K20").Value = 4827.8887
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = -4580.8887
$ws.Range("N20").Value = -6494
$ws.Range("H26").Value = 24980.334
$ws.Range("I26").Value = 9976.6
$ws.Range("K26").Value = 9976.6
$ws.Range("M26").Value = -9684.6
$ws.Range("H96").Value = 30000
$ws.Range("I96").Value = 30000
$ws.Range("K96").Value = 30000
$ws.Range("M96").Value = -27254
$ws.Range("H99").Value = 1928.4445
$ws.Range("I99").Value = 1669.5
$ws.Range("K99").Value = 1669.5
$ws.Range("M99").Value = -171.5
$ws.Range("H105").Value = 1043162.44
$ws.Range("I105").Value = 2287188.5
$ws.Range("K105").Value = 2287188.5
$ws.Range("M105").Value = -2285441.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21741340
$ws.Range("I31").Value = 25642788
$ws.Range("J31").Value = 4698.2856
$ws.Range("K31").Value = 25642788
$ws.Range("L31").Value = 4698.2856
$ws.Range("M31").Value = -25642493
$ws.Range("N31").Value = -5288.2856
$ws.Range("H34").Value = 21741340
$ws.Range("I34").Value = 25642788
$ws.Range("J34").Value = 4698.2856
$ws.Range("K34").Value = 25642788
$ws.Range("L34").Value = 4698.2856
$ws.Range("M34").Value = -25642586
$ws.Range("N34").Value = -5102.2856
$ws.Range("H99").Value = 8828.069
$ws.Range("I99").Value = 8355.647000000001
$ws.Range("J99").Value = 9497.333000000001
$ws.Range("K99").Value = 8355.647000000001
$ws.Range("L99").Value = 9497.333000000001
$ws.Range("M99").Value = -6857.647000000001
$ws.Range("N99").Value = -12493.333
$ws.Range("H126").Value = 8828.069
$ws.Range("I126").Value = 8355.647000000001
$ws.Range("J126").Value = 9497.333000000001
$ws.Range("K126").Value = 25066.941
$ws.Range("L126").Value = 28491.999
$ws.Range("M126").Value = -22596.941
$ws.Range("N126").Value = -33431.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2019265.6
$ws.Range("J32").Value = 2274082
$ws.Range("L32").Value = 6822246
$ws.Range("N32").Value = -6822812
$ws.Range("H38").Value = 54
$ws.Range("I38").Value = 40.285713
$ws.Range("K38").Value = 120.857139
$ws.Range("M38").Value = 226.142861
$ws.Range("H121").Value = 5450.3335
$ws.Range("J121").Value = 5838.9287
$ws.Range("L121").Value = 17516.7861
$ws.Range("N121").Value = -20136.7861
$ws.Range("H131").Value = 3887.1738
$ws.Range("I131").Value = 2393.182
$ws.Range("J131").Value = 5256.6665
$ws.Range("K131").Value = 7179.545999999999
$ws.Range("L131").Value = 15769.9995
$ws.Range("M131").Value = -2139.545999999999
$ws.Range("N131").Value = -25849.9995
$ws.Range("H132").Value = 1612.3636
$ws.Range("I132").Value = 666.1667
$ws.Range("J132").Value = 1967.1875
$ws.Range("K132").Value = 5995.5003
$ws.Range("L132").Value = 17704.6875
$ws.Range("M132").Value = -3465.5003
$ws.Range("N132").Value = -22764.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 8297.125
$ws.Range("J97").Value = 11321.363
$ws.Range("L97").Value = 11321.363
$ws.Range("N97").Value = -12313.363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2933
$ws.Range("I22").Value = 1950
$ws.Range("J22").Value = 3084.2307
$ws.Range("K22").Value = 1950
$ws.Range("L22").Value = 3084.2307
$ws.Range("M22").Value = -1655
$ws.Range("N22").Value = -3674.2307
$ws.Range("H27").Value = 2933
$ws.Range("I27").Value = 1950
$ws.Range("J27").Value = 3084.2307
$ws.Range("K27").Value = 1950
$ws.Range("L27").Value = 3084.2307
$ws.Range("M27").Value = -1843
$ws.Range("N27").Value = -3298.2307
$ws.Range("H40").Value = 6083.5
$ws.Range("I40").Value = 6083.5
$ws.Range("K40").Value = 6083.5
$ws.Range("M40").Value = -5947.5
$ws.Range("H46").Value = 1398.2
$ws.Range("J46").Value = 1998.3334
$ws.Range("L46").Value = 1998.3334
$ws.Range("N46").Value = -2374.3334
$ws.Range("H55").Value = 1001.5926
$ws.Range("I55").Value = 641.55554
$ws.Range("J55").Value = 1181.6111
$ws.Range("K55").Value = 641.55554
$ws.Range("L55").Value = 1181.6111
$ws.Range("M55").Value = -468.55554
$ws.Range("N55").Value = -1527.6111
$ws.Range("H82").Value = 2203.2856
$ws.Range("J82").Value = 3653.0833
$ws.Range("L82").Value = 3653.0833
$ws.Range("N82").Value = -4375.0833
$ws.Range("H85").Value = 2203.2856
$ws.Range("J85").Value = 3653.0833
$ws.Range("L85").Value = 3653.0833
$ws.Range("N85").Value = -6149.0833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 2340.6667
$ws.Range("I20").Value = 1005.5
$ws.Range("J20").Value = 5011
$ws.Range("K20").Value = 1005.5
$ws.Range("L20").Value = 5011
$ws.Range("M20").Value = -765.5
$ws.Range("N20").Value = -5491
$ws.Range("H45").Value = 19075.5
$ws.Range("J45").Value = 13371
$ws.Range("L45").Value = 13371
$ws.Range("N45").Value = -14353
$ws.Range("H132").Value = 563075.4
$ws.Range("J132").Value = 2501800
$ws.Range("L132").Value = 7505400
$ws.Range("N132").Value = -7510460
